$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# --- Fix the "fortrængning" calculation mistake for all biogas scenarios ---
# Row 2: biogas / svin
$ws.Range("C2").Value = 1036.72002250176
$ws.Range("D2").Value = 82.9459786893716

# Row 3: biogas / kvæg
$ws.Range("C3").Value = 708.75764551432405
$ws.Range("D3").Value = 65.954871804325506

# Row 11: køling_biogas / svin
$ws.Range("C11").Value = 1134.1193080038599
$ws.Range("D11").Value = 93.074749339300894

# Row 19: ugentlig_biogas / svin
$ws.Range("C19").Value = 1226.5554600062001
$ws.Range("D19").Value = 110.482338950433

# Row 20: ugentlig_biogas / kvæg
$ws.Range("C20").Value = 469.68380389112099
$ws.Range("D20").Value = 93.221355968363

# --- Resize / refit the columns so the table content is fully visible ---
$ws.Columns.Item(1).ColumnWidth = 33.59
$ws.Columns.Item(2).ColumnWidth = 14.59
$ws.Columns.Item(3).ColumnWidth = 28.09
$ws.Columns.Item(4).ColumnWidth = 33.42
$ws.Columns.Item(5).ColumnWidth = 12.75
$ws.Columns.Item(6).ColumnWidth = 12.92
$ws.Columns.Item(7).ColumnWidth = 11.09
$ws.Columns.Item(8).ColumnWidth = 11.42
$ws.Columns.Item(9).ColumnWidth = 11.09
$ws.Columns.Item(10).ColumnWidth = 11.75

# --- Update the active selection on the sheet ---
$ws.Range("M8").Select() | Out-Null
